$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Copy formatting (fill/border) from already "done" cells onto the cells
#     whose progress state changes, mirroring how this was done by hand with
#     the Format Painter, then update their text. ---

# D4, D6, D11 take on the same green "X" look as D5
$ws.Range("D5").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4122) | Out-Null
$ws.Range("D6").PasteSpecial(-4122) | Out-Null
$ws.Range("D11").PasteSpecial(-4122) | Out-Null

# E4, E6, E11, G13 take on the same green "X" look as E5
$ws.Range("E5").Copy() | Out-Null
$ws.Range("E4").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").PasteSpecial(-4122) | Out-Null
$ws.Range("E11").PasteSpecial(-4122) | Out-Null
$ws.Range("G13").PasteSpecial(-4122) | Out-Null

# H7 becomes a "problème" cell like H10
$ws.Range("H10").Copy() | Out-Null
$ws.Range("H7").PasteSpecial(-4122) | Out-Null

# G19 loses its own special "Commencé" box and matches the rest of row 19
$ws.Range("E19").Copy() | Out-Null
$ws.Range("G19").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Update cell contents to reflect the new progress state ---
$ws.Range("D4").Value2 = "X"
$ws.Range("E4").Value2 = "X"
$ws.Range("D6").Value2 = "X"
$ws.Range("E6").Value2 = "X"
$ws.Range("H7").Value2 = "problème"
$ws.Range("D11").Value2 = "X"
$ws.Range("E11").Value2 = "X"
$ws.Range("G13").Value2 = "X"
$ws.Range("G18").Value2 = "presque fini"
$ws.Range("G19").Value2 = "X"

# --- Restore the last active selection ---
$ws.Range("I26").Select() | Out-Null
